# Insert two new rows of data at row 921, pushing the existing rows
# (921-976) down to (923-978), matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 921 (existing rows shift down).
$ws.Rows.Item(921).Insert()
$ws.Rows.Item(921).Insert()

# Populate new row 921.
$ws.Cells.Item(921, 1).Value = 6
$ws.Cells.Item(921, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(921, 3).Value = "Metropolitana"
$ws.Cells.Item(921, 4).Value = 44826
$ws.Cells.Item(921, 5).Value = 13
$ws.Cells.Item(921, 6).Value = 100112031
$ws.Cells.Item(921, 7).Value = "Poroto verde"
$ws.Cells.Item(921, 8).Value = "Magnum"
$ws.Cells.Item(921, 9).Value = "Primera"
$ws.Cells.Item(921, 10).Value = 120
$ws.Cells.Item(921, 11).Value = 26000
$ws.Cells.Item(921, 12).Value = 27000
$ws.Cells.Item(921, 13).Value = 26583
$ws.Cells.Item(921, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(921, 15).Value = "Perú"
$ws.Cells.Item(921, 16).Value = 1063
$ws.Cells.Item(921, 17).Value = 25
$ws.Cells.Item(921, 18).Value = "Hortaliza"

# Populate new row 922.
$ws.Cells.Item(922, 1).Value = 6
$ws.Cells.Item(922, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(922, 3).Value = "Metropolitana"
$ws.Cells.Item(922, 4).Value = 44826
$ws.Cells.Item(922, 5).Value = 13
$ws.Cells.Item(922, 6).Value = 100112031
$ws.Cells.Item(922, 7).Value = "Poroto verde"
$ws.Cells.Item(922, 8).Value = "Sin especificar"
$ws.Cells.Item(922, 9).Value = "Primera"
$ws.Cells.Item(922, 10).Value = 35
$ws.Cells.Item(922, 11).Value = 48000
$ws.Cells.Item(922, 12).Value = 48000
$ws.Cells.Item(922, 13).Value = 48000
$ws.Cells.Item(922, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(922, 15).Value = "Perú"
$ws.Cells.Item(922, 16).Value = 1920
$ws.Cells.Item(922, 17).Value = 25
$ws.Cells.Item(922, 18).Value = "Hortaliza"

# Make sure the date cells keep the existing date style (column D, style index 2
# in the original file) in case Insert() didn't propagate it.
$ws.Cells.Item(921, 4).NumberFormat = $ws.Cells.Item(923, 4).NumberFormat
$ws.Cells.Item(922, 4).NumberFormat = $ws.Cells.Item(923, 4).NumberFormat
